$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of results (row 6) matching the existing data layout
$ws.Range("A6").Value = 42602.584131944444
$ws.Range("A6").NumberFormat = "m/d/yy h:mm"

$ws.Range("B6").Value = "Noun"

$ws.Range("C6").Value = 8516
$ws.Range("D6").Value = 3947
$ws.Range("E6").Value = 656
$ws.Range("F6").Value = 95
$ws.Range("G6").Value = 38
$ws.Range("H6").Value = 71
$ws.Range("I6").Value = 28
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 99
$ws.Range("M6").Value = 0
